$d = $word.ActiveDocument

# Locate the paragraph that must be kept: "LOQ4031: Química Geral I (Requisito)"
$anchor = $d.Content
$anchor.Find.Execute("LOQ4031: Química Geral I (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchorParaIndex = $anchor.Paragraphs.First.Index

# Locate the last paragraph that must be removed: the copyright/footer line
$footer = $d.Content
$footer.Find.Execute("© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$footerParaIndex = $footer.Paragraphs.First.Index

# Build a range spanning everything from just after the anchor paragraph
# through (and including) the footer paragraph's own paragraph mark, so the
# whole block of paragraphs (trailing blank line, "Ver no Jupiter..." line,
# and the copyright line) is removed in one go.
$startPara = $d.Paragraphs.Item($anchorParaIndex + 1)
$endPara = $d.Paragraphs.Item($footerParaIndex)

$delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$delRange.Delete()
